# Historical RPK data entry: the RPKs (mils) column (K) for rows 73-83 (years
# 2000-2010) was stored as German-locale decimal-comma text (e.g.
# "3201366,12411332"). Replace each with the plain (truncated) numeric value,
# matching the numeric style already used for the later years (rows 84+).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K73").Value = 3201366
$ws.Range("K74").Value = 3108526
$ws.Range("K75").Value = 3124069
$ws.Range("K76").Value = 3180302
$ws.Range("K77").Value = 3628725
$ws.Range("K78").Value = 3919023
$ws.Range("K79").Value = 4170556
$ws.Range("K80").Value = 4513095
$ws.Range("K81").Value = 4608466
$ws.Range("K82").Value = 4561413
$ws.Range("K83").Value = 4930250

# Mirror the author's final on-screen selection (scrolled down to row 84,
# cell K84 active) when they saved the workbook.
$ws.Range("K84").Select()
